$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 218-222 with revised figures
# Row 218
$ws.Range("B218").Value = 210862
$ws.Range("C218").Value = 70730
$ws.Range("D218").Value = 36657
$ws.Range("E218").Value = 1297
$ws.Range("F218").Value = 5001
$ws.Range("G218").Value = 27774
$ws.Range("H218").Value = 140132
$ws.Range("I218").Value = 19142
$ws.Range("J218").Value = 6902
$ws.Range("K218").Value = 114087
$ws.Range("L218").Value = 20172
$ws.Range("M218").Value = 3022
$ws.Range("N218").Value = 159
$ws.Range("O218").Value = 0
$ws.Range("P218").Value = 694
$ws.Range("Q218").Value = 2169
$ws.Range("R218").Value = 17150
$ws.Range("S218").Value = 6377
$ws.Range("T218").Value = 141
$ws.Range("U218").Value = 10632
$ws.Range("V218").Value = 190690
$ws.Range("W218").Value = 67708
$ws.Range("X218").Value = 36498
$ws.Range("Y218").Value = 1297
$ws.Range("Z218").Value = 4307
$ws.Range("AA218").Value = 25606
$ws.Range("AB218").Value = 122982
$ws.Range("AC218").Value = 12765
$ws.Range("AD218").Value = 6761
$ws.Range("AE218").Value = 103456

# Row 219
$ws.Range("B219").Value = 209836
$ws.Range("C219").Value = 69018
$ws.Range("D219").Value = 35564
$ws.Range("E219").Value = 1655
$ws.Range("F219").Value = 4874
$ws.Range("G219").Value = 26924
$ws.Range("H219").Value = 140818
$ws.Range("I219").Value = 19114
$ws.Range("J219").Value = 7032
$ws.Range("K219").Value = 114671
$ws.Range("L219").Value = 21168
$ws.Range("M219").Value = 3328
$ws.Range("N219").Value = 163
$ws.Range("O219").Value = 372
$ws.Range("P219").Value = 583
$ws.Range("Q219").Value = 2210
$ws.Range("R219").Value = 17840
$ws.Range("S219").Value = 6511
$ws.Range("T219").Value = 152
$ws.Range("U219").Value = 11176
$ws.Range("V219").Value = 188668
$ws.Range("W219").Value = 65690
$ws.Range("X219").Value = 35402
$ws.Range("Y219").Value = 1283
$ws.Range("Z219").Value = 4291
$ws.Range("AA219").Value = 24714
$ws.Range("AB219").Value = 122978
$ws.Range("AC219").Value = 12603
$ws.Range("AD219").Value = 6880
$ws.Range("AE219").Value = 103495

# Row 220
$ws.Range("B220").Value = 209510
$ws.Range("C220").Value = 66867
$ws.Range("D220").Value = 33754
$ws.Range("E220").Value = 1502
$ws.Range("F220").Value = 4855
$ws.Range("G220").Value = 26756
$ws.Range("H220").Value = 142642
$ws.Range("I220").Value = 19912
$ws.Range("J220").Value = 6875
$ws.Range("K220").Value = 115856
$ws.Range("L220").Value = 23121
$ws.Range("M220").Value = 3563
$ws.Range("N220").Value = 280
$ws.Range("O220").Value = 242
$ws.Range("P220").Value = 762
$ws.Range("Q220").Value = 2278
$ws.Range("R220").Value = 19558
$ws.Range("S220").Value = 7469
$ws.Range("T220").Value = 155
$ws.Range("U220").Value = 11935
$ws.Range("V220").Value = 186389
$ws.Range("W220").Value = 63305
$ws.Range("X220").Value = 33474
$ws.Range("Y220").Value = 1260
$ws.Range("Z220").Value = 4093
$ws.Range("AA220").Value = 24477
$ws.Range("AB220").Value = 123084
$ws.Range("AC220").Value = 12443
$ws.Range("AD220").Value = 6720
$ws.Range("AE220").Value = 103921

# Row 221
$ws.Range("B221").Value = 213397
$ws.Range("C221").Value = 68873
$ws.Range("D221").Value = 35862
$ws.Range("E221").Value = 1275
$ws.Range("F221").Value = 4835
$ws.Range("G221").Value = 26901
$ws.Range("H221").Value = 144524
$ws.Range("I221").Value = 19759
$ws.Range("J221").Value = 6977
$ws.Range("K221").Value = 117788
$ws.Range("L221").Value = 22599
$ws.Range("M221").Value = 3267
$ws.Range("N221").Value = 213
$ws.Range("O221").Value = 0
$ws.Range("P221").Value = 782
$ws.Range("Q221").Value = 2272
$ws.Range("R221").Value = 19332
$ws.Range("S221").Value = 7281
$ws.Range("T221").Value = 155
$ws.Range("U221").Value = 11896
$ws.Range("V221").Value = 190797
$ws.Range("W221").Value = 65606
$ws.Range("X221").Value = 35650
$ws.Range("Y221").Value = 1275
$ws.Range("Z221").Value = 4052
$ws.Range("AA221").Value = 24629
$ws.Range("AB221").Value = 125192
$ws.Range("AC221").Value = 12478
$ws.Range("AD221").Value = 6822
$ws.Range("AE221").Value = 105891

# Row 222
$ws.Range("B222").Value = 214870
$ws.Range("C222").Value = 69923
$ws.Range("D222").Value = 37594
$ws.Range("E222").Value = 1290
$ws.Range("F222").Value = 4716
$ws.Range("G222").Value = 26323
$ws.Range("H222").Value = 144946
$ws.Range("I222").Value = 20422
$ws.Range("J222").Value = 6928
$ws.Range("K222").Value = 117596
$ws.Range("L222").Value = 23465
$ws.Range("M222").Value = 3129
$ws.Range("N222").Value = 143
$ws.Range("O222").Value = 0
$ws.Range("P222").Value = 666
$ws.Range("Q222").Value = 2320
$ws.Range("R222").Value = 20335
$ws.Range("S222").Value = 8386
$ws.Range("T222").Value = 160
$ws.Range("U222").Value = 11789
$ws.Range("V222").Value = 191405
$ws.Range("W222").Value = 66794
$ws.Range("X222").Value = 37451
$ws.Range("Y222").Value = 1290
$ws.Range("Z222").Value = 4050
$ws.Range("AA222").Value = 24003
$ws.Range("AB222").Value = 124611
$ws.Range("AC222").Value = 12036
$ws.Range("AD222").Value = 6768
$ws.Range("AE222").Value = 105807

# New row 223 (01-06-2021)
$ws.Range("A223").Formula = '="01-06-2021"'
$ws.Range("A223").Copy()
$ws.Range("A223").PasteSpecial(-4163)

$ws.Range("B223").Value = 213283
$ws.Range("C223").Value = 70287
$ws.Range("D223").Value = 37860
$ws.Range("E223").Value = 1265
$ws.Range("F223").Value = 4616
$ws.Range("G223").Value = 26546
$ws.Range("H223").Value = 142996
$ws.Range("I223").Value = 19660
$ws.Range("J223").Value = 6850
$ws.Range("K223").Value = 116486
$ws.Range("L223").Value = 22318
$ws.Range("M223").Value = 2814
$ws.Range("N223").Value = 142
$ws.Range("O223").Value = 0
$ws.Range("P223").Value = 603
$ws.Range("Q223").Value = 2069
$ws.Range("R223").Value = 19504
$ws.Range("S223").Value = 7843
$ws.Range("T223").Value = 178
$ws.Range("U223").Value = 11482
$ws.Range("V223").Value = 190965
$ws.Range("W223").Value = 67473
$ws.Range("X223").Value = 37718
$ws.Range("Y223").Value = 1265
$ws.Range("Z223").Value = 4013
$ws.Range("AA223").Value = 24476
$ws.Range("AB223").Value = 123492
$ws.Range("AC223").Value = 11817
$ws.Range("AD223").Value = 6672
$ws.Range("AE223").Value = 105003
